$p = $ppt.ActivePresentation

# 1. Merge the three runs in the "For a / more accurate study / , we'd need:" paragraph
#    (slide 5, shape 4 "TextBox 13") into a single run. Using Characters() on the
#    whole paragraph span (rather than setting .Text on the Paragraphs() range
#    directly) forces a clean single-run replacement instead of a diff-preserving
#    run-boundary split.
$s5 = $p.Slides.Item(5)
$tb = $s5.Shapes.Item(4)
$tr = $tb.TextFrame.TextRange
$para = $tr.Paragraphs(9, 1)
$whole = $tr.Characters($para.Start, $para.Text.Length)
$whole.Text = "For a more accurate study, we’d need:"

# 2. Delete slide 6 ("TROUBLES 'N' TURMOIL") entirely.
$p.Slides.Item(6).Delete()
